$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9656778573989868
$ws.Range("B1").Value = 2.758779525756836
$ws.Range("C1").Value = 5.52840518951416
$ws.Range("D1").Value = 2.107406616210938
$ws.Range("E1").Value = 1.19182562828064
